# Trade #37 closed at 2026-02-16 21:28:53 - momentum DOWN +0.000%
# Append a new OPEN trade row (row 9) to the "momentum" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("momentum")

$row = 9

# A: Trade #
$ws.Cells.Item($row, 1).Value = 37

# B: Date - force text so "2026-02-16" isn't auto-parsed into a date serial
# (matches how the existing rows store it as literal text), then drop the
# temporary number format so the cell keeps the sheet's default style.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).ClearFormats()

# C: Time (kept as plain text, like the other rows)
$ws.Cells.Item($row, 3).Value = "21:28:53"

# D: Strategy
$ws.Cells.Item($row, 4).Value = "momentum"

# E: Side
$ws.Cells.Item($row, 5).Value = "DOWN"

# F: Entry Price
$ws.Cells.Item($row, 6).Value = 68582.45

# G: Exit Price - empty (trade still open). A plain "" assignment leaves the
# cell completely blank (no <c> element at all), so use the quote-prefix
# trick to force an explicit empty-string cell, then strip the formatting
# it leaves behind so the cell stays on the default style.
$ws.Cells.Item($row, 7).Value = "'"
$ws.Cells.Item($row, 7).ClearFormats()

# H: Status
$ws.Cells.Item($row, 8).Value = "OPEN"

# I: P&L %
$ws.Cells.Item($row, 9).Value = 0

# J: P&L $
$ws.Cells.Item($row, 10).Value = 0

# K: Confidence
$ws.Cells.Item($row, 11).Value = 0.9

# L: Entry Reason
$ws.Cells.Item($row, 12).Value = "Downward momentum: -0.488% over 10 samples"

# M: Exit Reason - empty (trade still open), same empty-cell trick as G.
$ws.Cells.Item($row, 13).Value = "'"
$ws.Cells.Item($row, 13).ClearFormats()

# N: Duration (min)
$ws.Cells.Item($row, 14).Value = 0
